$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (lower-case column names)
$ws.Range("A1").Value = "cnae"
$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "lp_exposure"

# Update LP_Exposure values for the specified rows
$updates = @{
    2  = "high"
    4  = "high"
    6  = "median"
    10 = "median"
    11 = "median"
    12 = "median"
    13 = "median"
    15 = "median"
    16 = "median"
    17 = "median"
    18 = "median"
    19 = "median"
    20 = "median"
    21 = "median"
    22 = "median"
    23 = "median"
    24 = "median"
    25 = "median"
    26 = "median"
    27 = "median"
    28 = "median"
    34 = "median"
    37 = "high"
    38 = "high"
    39 = "high"
    40 = "high"
    41 = "high"
    52 = "low"
    53 = "low"
    54 = "low"
    55 = "low"
    56 = "low"
    58 = "low"
    59 = "low"
    60 = "low"
    61 = "low"
    62 = "low"
    63 = "low"
    65 = "low"
    66 = "low"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
